$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin")
$ws.Range("A9").Value = "Standards file"
Write-Host "done"
